$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12 (Leve Item ID 5515)
$ws.Cells.Item(12, 8).Value = 289.46667
$ws.Cells.Item(12, 9).Value = 164.2
$ws.Cells.Item(12, 11).Value = 164.2
$ws.Cells.Item(12, 13).Value = 5.800000000000011

# Row 41 (Leve Item ID 5478)
$ws.Cells.Item(41, 8).Value = 279.5263
$ws.Cells.Item(41, 9).Value = 274.07144
$ws.Cells.Item(41, 10).Value = 294.8
$ws.Cells.Item(41, 11).Value = 274.07144
$ws.Cells.Item(41, 12).Value = 294.8
$ws.Cells.Item(41, 13).Value = 165.92856
$ws.Cells.Item(41, 14).Value = -1174.8

# Row 116 (Leve Item ID 27778)
$ws.Cells.Item(116, 8).Value = 4502.1665
$ws.Cells.Item(116, 9).Value = 5054.4443
$ws.Cells.Item(116, 11).Value = 5054.4443
$ws.Cells.Item(116, 13).Value = -1612.4443

# Row 137 (Leve Item ID 44013)
$ws.Cells.Item(137, 8).Value = 2159.1333
$ws.Cells.Item(137, 9).Value = 2168.2307
$ws.Cells.Item(137, 11).Value = 6504.6921
$ws.Cells.Item(137, 13).Value = -3954.6921

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Cells.Item(32, 8).Value = 6636.363
$ws.Cells.Item(32, 9).Value = 4672.1606
$ws.Cells.Item(32, 11).Value = 4672.1606
$ws.Cells.Item(32, 13).Value = -4385.1606

# Row 63 (Leve Item ID 12528)
$ws.Cells.Item(63, 8).Value = 142860880
$ws.Cells.Item(63, 9).Value = 200004500
$ws.Cells.Item(63, 10).Value = 1850
$ws.Cells.Item(63, 11).Value = 200004500
$ws.Cells.Item(63, 12).Value = 1850
$ws.Cells.Item(63, 13).Value = -200003814
$ws.Cells.Item(63, 14).Value = -3222

# Row 66 (Leve Item ID 12528)
$ws.Cells.Item(66, 8).Value = 142860880
$ws.Cells.Item(66, 9).Value = 200004500
$ws.Cells.Item(66, 10).Value = 1850
$ws.Cells.Item(66, 11).Value = 1000022500
$ws.Cells.Item(66, 12).Value = 9250
$ws.Cells.Item(66, 13).Value = -1000019068
$ws.Cells.Item(66, 14).Value = -16114

# Row 97 (Leve Item ID 19941)
$ws.Cells.Item(97, 8).Value = 983
$ws.Cells.Item(97, 9).Value = 701.2857
$ws.Cells.Item(97, 10).Value = 1640.3334
$ws.Cells.Item(97, 11).Value = 701.2857
$ws.Cells.Item(97, 12).Value = 1640.3334
$ws.Cells.Item(97, 13).Value = -205.2857
$ws.Cells.Item(97, 14).Value = -2632.3334

# Row 122 (Leve Item ID 36168)
$ws.Cells.Item(122, 8).Value = 828691.1
$ws.Cells.Item(122, 9).Value = 856274.2
$ws.Cells.Item(122, 10).Value = 1200
$ws.Cells.Item(122, 11).Value = 2568822.6
$ws.Cells.Item(122, 12).Value = 3600
$ws.Cells.Item(122, 13).Value = -2566372.6
$ws.Cells.Item(122, 14).Value = -8500

$ws = $wb.Worksheets.Item("BSM")
# Row 64 (Leve Item ID 14184)
$ws.Cells.Item(64, 8).Value = 226.75
$ws.Cells.Item(64, 9).Value = 169.625
$ws.Cells.Item(64, 10).Value = 264.83334
$ws.Cells.Item(64, 11).Value = 169.625
$ws.Cells.Item(64, 12).Value = 264.83334
$ws.Cells.Item(64, 13).Value = 55.375
$ws.Cells.Item(64, 14).Value = -714.83334

# Row 67 (Leve Item ID 14184)
$ws.Cells.Item(67, 8).Value = 226.75
$ws.Cells.Item(67, 9).Value = 169.625
$ws.Cells.Item(67, 10).Value = 264.83334
$ws.Cells.Item(67, 11).Value = 169.625
$ws.Cells.Item(67, 12).Value = 264.83334
$ws.Cells.Item(67, 13).Value = 610.375
$ws.Cells.Item(67, 14).Value = -1824.83334

# Row 75 (Leve Item ID 11872)
$ws.Cells.Item(75, 8).Value = 12465.6
$ws.Cells.Item(75, 9).Value = 8809.333000000001
$ws.Cells.Item(75, 10).Value = 17950
$ws.Cells.Item(75, 11).Value = 8809.333000000001
$ws.Cells.Item(75, 12).Value = 17950
$ws.Cells.Item(75, 13).Value = -7873.333000000001
$ws.Cells.Item(75, 14).Value = -19822

# Row 78 (Leve Item ID 11872)
$ws.Cells.Item(78, 8).Value = 12465.6
$ws.Cells.Item(78, 9).Value = 8809.333000000001
$ws.Cells.Item(78, 10).Value = 17950
$ws.Cells.Item(78, 11).Value = 26427.999
$ws.Cells.Item(78, 12).Value = 53850
$ws.Cells.Item(78, 13).Value = -21747.999
$ws.Cells.Item(78, 14).Value = -63210

# Row 107 (Leve Item ID 27706)
$ws.Cells.Item(107, 8).Value = 1008.96155
$ws.Cells.Item(107, 9).Value = 856.8889
$ws.Cells.Item(107, 10).Value = 1351.125
$ws.Cells.Item(107, 11).Value = 856.8889
$ws.Cells.Item(107, 12).Value = 1351.125
$ws.Cells.Item(107, 13).Value = 1063.1111
$ws.Cells.Item(107, 14).Value = -5191.125

# Row 134 (Leve Item ID 43998)
$ws.Cells.Item(134, 8).Value = 2730.36
$ws.Cells.Item(134, 9).Value = 2451.6667
$ws.Cells.Item(134, 10).Value = 2887.125
$ws.Cells.Item(134, 11).Value = 7355.000100000001
$ws.Cells.Item(134, 12).Value = 8661.375
$ws.Cells.Item(134, 13).Value = -4820.000100000001
$ws.Cells.Item(134, 14).Value = -13731.375

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (Leve Item ID 5367)
$ws.Cells.Item(22, 8).Value = 674.25
$ws.Cells.Item(22, 9).Value = 674.25
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 674.25
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -324.25
$ws.Cells.Item(22, 14).Value = $null

# Row 31 (Leve Item ID 44023)
$ws.Cells.Item(31, 8).Value = 3762.2957
$ws.Cells.Item(31, 9).Value = 1692.8536
$ws.Cells.Item(31, 10).Value = 6590.533
$ws.Cells.Item(31, 11).Value = 1692.8536
$ws.Cells.Item(31, 12).Value = 6590.533
$ws.Cells.Item(31, 13).Value = -1397.8536
$ws.Cells.Item(31, 14).Value = -7180.533

# Row 34 (Leve Item ID 44023)
$ws.Cells.Item(34, 8).Value = 3762.2957
$ws.Cells.Item(34, 9).Value = 1692.8536
$ws.Cells.Item(34, 10).Value = 6590.533
$ws.Cells.Item(34, 11).Value = 1692.8536
$ws.Cells.Item(34, 12).Value = 6590.533
$ws.Cells.Item(34, 13).Value = -1490.8536
$ws.Cells.Item(34, 14).Value = -6994.533

# Row 86 (Leve Item ID 12584)
$ws.Cells.Item(86, 8).Value = 3096.75
$ws.Cells.Item(86, 9).Value = 3190
$ws.Cells.Item(86, 10).Value = 3003.5
$ws.Cells.Item(86, 11).Value = 3190
$ws.Cells.Item(86, 12).Value = 3003.5
$ws.Cells.Item(86, 13).Value = -2067
$ws.Cells.Item(86, 14).Value = -5249.5

# Row 89 (Leve Item ID 12584)
$ws.Cells.Item(89, 8).Value = 3096.75
$ws.Cells.Item(89, 9).Value = 3190
$ws.Cells.Item(89, 10).Value = 3003.5
$ws.Cells.Item(89, 11).Value = 15950
$ws.Cells.Item(89, 12).Value = 15017.5
$ws.Cells.Item(89, 13).Value = -10334
$ws.Cells.Item(89, 14).Value = -26249.5

# Row 99 (Leve Item ID 36198)
$ws.Cells.Item(99, 8).Value = 6254.5454
$ws.Cells.Item(99, 10).Value = 1480
$ws.Cells.Item(99, 12).Value = 1480
$ws.Cells.Item(99, 14).Value = -4476

# Row 105 (Leve Item ID 19928)
$ws.Cells.Item(105, 8).Value = 1608.1364
$ws.Cells.Item(105, 9).Value = 1711.25
$ws.Cells.Item(105, 11).Value = 1711.25
$ws.Cells.Item(105, 13).Value = 35.75

# Row 122 (Leve Item ID 36196)
$ws.Cells.Item(122, 8).Value = 2528810.2
$ws.Cells.Item(122, 9).Value = 6945214
$ws.Cells.Item(122, 10).Value = 5150.857
$ws.Cells.Item(122, 11).Value = 20835642
$ws.Cells.Item(122, 12).Value = 15452.571
$ws.Cells.Item(122, 13).Value = -20833192
$ws.Cells.Item(122, 14).Value = -20352.571

# Row 126 (Leve Item ID 36198)
$ws.Cells.Item(126, 8).Value = 6254.5454
$ws.Cells.Item(126, 10).Value = 1480
$ws.Cells.Item(126, 12).Value = 4440
$ws.Cells.Item(126, 14).Value = -9380

$ws = $wb.Worksheets.Item("CUL")
# Row 32 (Leve Item ID 4731)
$ws.Cells.Item(32, 8).Value = 3625
$ws.Cells.Item(32, 10).Value = 3625
$ws.Cells.Item(32, 12).Value = 10875
$ws.Cells.Item(32, 14).Value = -11441

# Row 68 (Leve Item ID 12895)
$ws.Cells.Item(68, 8).Value = 1121.5714
$ws.Cells.Item(68, 9).Value = 800.6667
$ws.Cells.Item(68, 10).Value = 1209.091
$ws.Cells.Item(68, 11).Value = 2402.0001
$ws.Cells.Item(68, 12).Value = 3627.273
$ws.Cells.Item(68, 13).Value = -1591.0001
$ws.Cells.Item(68, 14).Value = -5249.272999999999

# Row 71 (Leve Item ID 12895)
$ws.Cells.Item(71, 8).Value = 1121.5714
$ws.Cells.Item(71, 9).Value = 800.6667
$ws.Cells.Item(71, 10).Value = 1209.091
$ws.Cells.Item(71, 11).Value = 7206.0003
$ws.Cells.Item(71, 12).Value = 10881.819
$ws.Cells.Item(71, 13).Value = -3150.0003
$ws.Cells.Item(71, 14).Value = -18993.819

# Row 86 (Leve Item ID 12892)
$ws.Cells.Item(86, 8).Value = 1255.25
$ws.Cells.Item(86, 9).Value = 1767.1666
$ws.Cells.Item(86, 10).Value = 743.3333
$ws.Cells.Item(86, 11).Value = 5301.4998
$ws.Cells.Item(86, 12).Value = 2229.9999
$ws.Cells.Item(86, 13).Value = -4115.4998
$ws.Cells.Item(86, 14).Value = -4601.9999

# Row 89 (Leve Item ID 12892)
$ws.Cells.Item(89, 8).Value = 1255.25
$ws.Cells.Item(89, 9).Value = 1767.1666
$ws.Cells.Item(89, 10).Value = 743.3333
$ws.Cells.Item(89, 11).Value = 15904.4994
$ws.Cells.Item(89, 12).Value = 6689.9997
$ws.Cells.Item(89, 13).Value = -9976.499400000001
$ws.Cells.Item(89, 14).Value = -18545.9997

# Row 121 (Leve Item ID 27878)
$ws.Cells.Item(121, 8).Value = 17433.2
$ws.Cells.Item(121, 10).Value = 34302.6
$ws.Cells.Item(121, 12).Value = 102907.8
$ws.Cells.Item(121, 14).Value = -105527.8

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (Leve Item ID 36248)
$ws.Cells.Item(40, 8).Value = 166668660
$ws.Cells.Item(40, 9).Value = 166668660
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 166668660
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).Value = -166668524
$ws.Cells.Item(40, 14).Value = $null

$ws = $wb.Worksheets.Item("WVR")
# Row 113 (Leve Item ID 27752)
$ws.Cells.Item(113, 8).Value = 877.8148
$ws.Cells.Item(113, 9).Value = 780.8077
$ws.Cells.Item(113, 11).Value = 2342.4231
$ws.Cells.Item(113, 13).Value = -172.4231
